$wb = $excel.ActiveWorkbook

# --- Sheet "TDES": remove the row containing "codecutil" ---
$tdes = $wb.Worksheets.Item("TDES")
$tdes.Rows.Item(5).Delete()

# Hyperlinks do not auto-shift with the row delete in this engine, so
# rebuild them pointing at their new (shifted-up) cells, then strip the
# auto-applied hyperlink style/value changes so the cells stay exactly
# as before (same text, same formatting).
$tdes.Hyperlinks.Delete()

$h6 = $tdes.Hyperlinks.Add($tdes.Range("A6"), "http://localhost:8087/tdes/console/getTemplateList", "", "Click to edit template")
$h6.TextToDisplay = "http://localhost:8087/tdes/console/getTemplateList"
$tdes.Range("A6").Style = "Normal"
$tdes.Range("A6").Value = "Test Document Fill"

$h7 = $tdes.Hyperlinks.Add($tdes.Range("A7"), "http://localhost:8087/tdes/console/getTemplateList", "", "Click to edit template")
$h7.TextToDisplay = "http://localhost:8087/tdes/console/getTemplateList"
$tdes.Range("A7").Style = "Normal"
$tdes.Range("A7").Value = "xmldocfill"

# --- Sheet "CV": remove the rows containing "PDFTools",
#     "Test - TempFile - SMBV2FileTransferUtils" and "TestSuite - TempFile" ---
$cv = $wb.Worksheets.Item("CV")
$cv.Rows.Item(38).Delete()
$cv.Rows.Item(34).Delete()
$cv.Rows.Item(15).Delete()

# --- Selections / active sheet ---
$tdes.Range("A5:XFD5").Select()
$cv.Range("A20").Select()
$cv.Activate()
